$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.826.12'
$ws.Range('D3').Value = '2.089.17'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.98'
$ws.Range('E5').Value = '  -0.52%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.68'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0785'
$ws.Range('E10').Value = '  -0.70%  '
$ws.Range('E11').Value = '  +2.76%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.19'
$ws.Range('E12').Value = '  +2.79%  '
$ws.Range('D13').Value = '2.397.68'
$ws.Range('E13').Value = '  +0.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.39'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.782'
$ws.Range('E15').Value = '  +1.27%  '
$ws.Range('D17').Value = '2.095.69'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('D18').Value = '37.787.06'
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('E19').Value = '  -0.92%  '
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').Value = '0.0₃0838'
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.41'
$ws.Range('E22').Value = '  +0.80%  '
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('E25').Value = '  -1.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.85'
$ws.Range('E26').Value = '  +9.30%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '171.95'
$ws.Range('E27').Value = '  +1.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.136'
$ws.Range('E28').Value = '  -2.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.53'
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('E32').Value = '  +0.46%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  -0.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.51'
$ws.Range('E35').Value = '  +0.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.82'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.40'
$ws.Range('E37').Value = '  -1.88%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.38'
$ws.Range('E39').Value = '  +0.12%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0236'
$ws.Range('E40').Value = '  +9.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '102.82'
$ws.Range('E41').Value = '  +4.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0976'
$ws.Range('E42').Value = '  -2.03%  '
$ws.Range('E43').Value = '  +0.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.86'
$ws.Range('E44').Value = '  +5.06%  '
$ws.Range('D45').Value = '1.456.80'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('E46').Value = '  -0.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.18'
$ws.Range('E47').Value = '  -4.54%  '
$ws.Range('E48').Value = '  -1.29%  '
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('E50').Value = '  -1.91%  '
$ws.Range('E51').Value = '  +0.01%  '
